# Insert a new weekly Chirimoya price record at row 349 of Sheet1.
# This pushes the existing rows 349-380 down to 350-381 (dimension
# grows from A1:T380 to A1:T381) and fills the freshly inserted
# row 349 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 349, shifting 349:380 -> 350:381.
$ws.Rows(349).Insert()

# Populate the new row 349 with the latest observation.
$ws.Cells.Item(349, 1).Value  = 6
$ws.Cells.Item(349, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(349, 3).Value  = "Metropolitana"
$ws.Cells.Item(349, 4).Value  = 45127
$ws.Cells.Item(349, 5).Value  = 13
$ws.Cells.Item(349, 6).Value  = "Fruta"
$ws.Cells.Item(349, 7).Value  = 100107
$ws.Cells.Item(349, 8).Value  = "Otros"
$ws.Cells.Item(349, 9).Value  = 100107002
$ws.Cells.Item(349, 10).Value = "Chirimoya"
$ws.Cells.Item(349, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(349, 12).Value = "Primera"
$ws.Cells.Item(349, 13).Value = 175
$ws.Cells.Item(349, 14).Value = 28000
$ws.Cells.Item(349, 15).Value = 28000
$ws.Cells.Item(349, 16).Value = 28000
$ws.Cells.Item(349, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(349, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(349, 19).Value = 2800
$ws.Cells.Item(349, 20).Value = 10
